$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B28: was stored as an inline string "3", should be a true number 3
$ws.Range("B28").Value = 3

# Add new row 29 with the additional annotation record
$ws.Range("A29").Value = "Ying Tang"
$ws.Range("B29").Value = "4"
$ws.Range("C29").Value = " appealing ,an important open problem"
$ws.Range("D29").Value = "APC"
$ws.Range("E29").Value = "MET"
$ws.Range("F29").Value = "f913699b-da49-47c6-8043-88c593733ae2"
$ws.Range("G29").Value = "BJyy3a0Ez_annotated.xlsx"
$ws.Range("H29").Value = "The idea of model-parallelism (as opposed to data parallelism) is appealing and an important open problem."
